# Daily attendance processing - 2025-12-22 19:49:37
#
# Reorders the comma-separated "Recorded By" names in column G so that the
# "System" entry and its adjacent collaborator entry (dnasr281@gmail.com or
# backup@backdoor.com) swap places, e.g.:
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System"         -> "System, backup@backdoor.com"
#   "system, backup@backdoor.com, System" -> "system, System, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7   # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $a = $parts[$parts.Count - 2]
    $b = $parts[$parts.Count - 1]

    $swap = $false
    if ($a -eq "System" -and ($b -eq "dnasr281@gmail.com" -or $b -eq "backup@backdoor.com")) {
        $swap = $true
    }
    elseif ($b -eq "System" -and ($a -eq "dnasr281@gmail.com" -or $a -eq "backup@backdoor.com")) {
        $swap = $true
    }

    if ($swap) {
        $parts[$parts.Count - 2] = $b
        $parts[$parts.Count - 1] = $a
        $newVal = [string]::Join(", ", $parts)
        $cell.Value = $newVal
    }
}
